# Update Slovakia Covid DailyStats worksheet (commit: "Updated: st 12. 08. 2021")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Corrections to existing AgTests (F) / AgPosit (G) values ---
$ws.Range("F397").Value = 107978

$ws.Range("F473").Value = 40174

$ws.Range("F480").Value = 33971

$ws.Range("F495").Value = 10367
$ws.Range("F496").Value = 8209
$ws.Range("F497").Value = 7695
$ws.Range("F498").Value = 9125
$ws.Range("F499").Value = 11341
$ws.Range("F500").Value = 7662
$ws.Range("F501").Value = 5687
$ws.Range("F502").Value = 10479
$ws.Range("F503").Value = 7408
$ws.Range("F504").Value = 7437
$ws.Range("F505").Value = 8474
$ws.Range("F506").Value = 10763
$ws.Range("F507").Value = 7149
$ws.Range("F508").Value = 5614
$ws.Range("F509").Value = 9494
$ws.Range("F510").Value = 7793
$ws.Range("F511").Value = 6754
$ws.Range("F512").Value = 8403
$ws.Range("F513").Value = 10283
$ws.Range("F514").Value = 6889
$ws.Range("F515").Value = 4960

$ws.Range("F516").Value = 9220
$ws.Range("G516").Value = 33

$ws.Range("F517").Value = 6670
$ws.Range("F518").Value = 6963

$ws.Range("F519").Value = 7763
$ws.Range("G519").Value = 20

$ws.Range("F520").Value = 9885
$ws.Range("G520").Value = 23

$ws.Range("F521").Value = 6159
$ws.Range("G521").Value = 21

$ws.Range("F522").Value = 4885
$ws.Range("G522").Value = 14

# --- Append three new rows of daily data (523-525) ---
$newRows = @(
    @{ Row = 523; A = 44417; B = 393160; C = 8193; D = 101; E = 12543; F = 9346; G = 18 },
    @{ Row = 524; A = 44418; B = 393228; C = 5739; D = 68;  E = 12544; F = 6948; G = 27 },
    @{ Row = 525; A = 44419; B = 393302; C = 6658; D = 74;  E = 12544; F = 5524; G = 18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
